$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsZhCn.Range("D4").Value = "2016-02-17 04:04:12"
$wsDeDe.Range("D4").Value = "2016-02-17 04:04:21"
